# Fix list level numbering: decrement the indent level of list-continuation
# paragraphs that were erroneously nested one level too deep.
#
# PowerPoint's TextRange.IndentLevel is 1-based (level 1 == OOXML a:pPr lvl="0"),
# so decrementing IndentLevel by 1 here corresponds to decrementing the OOXML
# lvl attribute by 1 in the diff.

$p = $ppt.ActivePresentation

# Slide 1, Content Placeholder: "Bullet item with inline code" (lvl 1 -> 0)
# and "with nested" (lvl 2 -> 1).
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(2).TextFrame.TextRange
$tr1.Paragraphs(2).IndentLevel = $tr1.Paragraphs(2).IndentLevel - 1
$tr1.Paragraphs(4).IndentLevel = $tr1.Paragraphs(4).IndentLevel - 1

# Slide 2, Content Placeholder: "Nested" (lvl 1 -> 0).
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Paragraphs(2).IndentLevel = $tr2.Paragraphs(2).IndentLevel - 1

# Slide 3, Content Placeholder: "A total alternative for head" (lvl 1 -> 0).
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$tr3.Paragraphs(1).IndentLevel = $tr3.Paragraphs(1).IndentLevel - 1
